# Generate Report for Handoff
# Update the "latest handoff" timestamps for files that are currently
# queued for handoff (status "Handback transform failed" or
# "Ready for handoff") to reflect a freshly generated handoff report.
# Files already "Handed back: in sync with en-US" or "In Translation"
# are left untouched.

$wb = $excel.ActiveWorkbook

$rowsToUpdate = @(4, 6, 7, 8, 9, 10)

# Overview sheet: column D holds "Latest Handoff Date"
$ws1 = $wb.Worksheets.Item("Overview")
foreach ($r in $rowsToUpdate) {
    $ws1.Cells.Item($r, 4).Value = "2016-03-20 05:29:20"
}

# zh-cn sheet: column E holds "Latest Handoff Datetime"
$ws2 = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rowsToUpdate) {
    $ws2.Cells.Item($r, 5).Value = "2016-03-20 05:29:11"
}

# de-de sheet: column E holds "Latest Handoff Datetime"
$ws3 = $wb.Worksheets.Item("de-de")
foreach ($r in $rowsToUpdate) {
    $ws3.Cells.Item($r, 5).Value = "2016-03-20 05:29:20"
}
